$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (style) from A1 into the full new header range P1:CJ1
# so that newly created cells beyond the old BA1 boundary get the same bold/border/center style.
$ws.Range("A1").Copy()
$ws.Range("P1:CJ1").PasteSpecial(-4122)

# Set header text for the new/shifted columns (P.. CJ). Columns A..O are unchanged.
$ws.Range("P1").Value = "companies.1.company [Link] [Company]"
$ws.Range("Q1").Value = "companies.2.company [Link] [Company]"
$ws.Range("R1").Value = "companies.3.company [Link] [Company]"
$ws.Range("S1").Value = "companies.4.company [Link] [Company]"
$ws.Range("T1").Value = "companies.5.company [Link] [Company]"
$ws.Range("U1").Value = "market_segment [Link] [Market Segment]"
$ws.Range("V1").Value = "industry [Link] [Industry Type]"
$ws.Range("W1").Value = "website [Data]"
$ws.Range("X1").Value = "language [Link] [Language]"
$ws.Range("Y1").Value = "customer_details [Text]"
$ws.Range("Z1").Value = "customer_primary_address [Link] [Address]"
$ws.Range("AA1").Value = "customer_primary_contact [Link] [Contact]"
$ws.Range("AB1").Value = "tax_id [Data]"
$ws.Range("AC1").Value = "tax_category [Link] [Tax Category]"
$ws.Range("AD1").Value = "tax_withholding_category [Link] [Tax Withholding Category]"
$ws.Range("AE1").Value = "payment_terms [Link] [Payment Terms Template]"
$ws.Range("AF1").Value = "credit_limits.1.company [Link] [Company]"
$ws.Range("AG1").Value = "credit_limits.1.credit_limit [Currency]"
$ws.Range("AH1").Value = "credit_limits.1.bypass_credit_limit_check [Check]"
$ws.Range("AI1").Value = "credit_limits.2.company [Link] [Company]"
$ws.Range("AJ1").Value = "credit_limits.2.credit_limit [Currency]"
$ws.Range("AK1").Value = "credit_limits.2.bypass_credit_limit_check [Check]"
$ws.Range("AL1").Value = "credit_limits.3.company [Link] [Company]"
$ws.Range("AM1").Value = "credit_limits.3.credit_limit [Currency]"
$ws.Range("AN1").Value = "credit_limits.3.bypass_credit_limit_check [Check]"
$ws.Range("AO1").Value = "credit_limits.4.company [Link] [Company]"
$ws.Range("AP1").Value = "credit_limits.4.credit_limit [Currency]"
$ws.Range("AQ1").Value = "credit_limits.4.bypass_credit_limit_check [Check]"
$ws.Range("AR1").Value = "credit_limits.5.company [Link] [Company]"
$ws.Range("AS1").Value = "credit_limits.5.credit_limit [Currency]"
$ws.Range("AT1").Value = "credit_limits.5.bypass_credit_limit_check [Check]"
$ws.Range("AU1").Value = "accounts.1.company [Link] [Company]"
$ws.Range("AV1").Value = "accounts.1.account [Link] [Account]"
$ws.Range("AW1").Value = "accounts.1.advance_account [Link] [Account]"
$ws.Range("AX1").Value = "accounts.2.company [Link] [Company]"
$ws.Range("AY1").Value = "accounts.2.account [Link] [Account]"
$ws.Range("AZ1").Value = "accounts.2.advance_account [Link] [Account]"
$ws.Range("BA1").Value = "accounts.3.company [Link] [Company]"
$ws.Range("BB1").Value = "accounts.3.account [Link] [Account]"
$ws.Range("BC1").Value = "accounts.3.advance_account [Link] [Account]"
$ws.Range("BD1").Value = "accounts.4.company [Link] [Company]"
$ws.Range("BE1").Value = "accounts.4.account [Link] [Account]"
$ws.Range("BF1").Value = "accounts.4.advance_account [Link] [Account]"
$ws.Range("BG1").Value = "accounts.5.company [Link] [Company]"
$ws.Range("BH1").Value = "accounts.5.account [Link] [Account]"
$ws.Range("BI1").Value = "accounts.5.advance_account [Link] [Account]"
$ws.Range("BJ1").Value = "loyalty_program [Link] [Loyalty Program]"
$ws.Range("BK1").Value = "sales_team.1.sales_person [Link] [Sales Person]"
$ws.Range("BL1").Value = "sales_team.1.allocated_percentage [Float]"
$ws.Range("BM1").Value = "sales_team.1.incentives [Currency]"
$ws.Range("BN1").Value = "sales_team.2.sales_person [Link] [Sales Person]"
$ws.Range("BO1").Value = "sales_team.2.allocated_percentage [Float]"
$ws.Range("BP1").Value = "sales_team.2.incentives [Currency]"
$ws.Range("BQ1").Value = "sales_team.3.sales_person [Link] [Sales Person]"
$ws.Range("BR1").Value = "sales_team.3.allocated_percentage [Float]"
$ws.Range("BS1").Value = "sales_team.3.incentives [Currency]"
$ws.Range("BT1").Value = "sales_team.4.sales_person [Link] [Sales Person]"
$ws.Range("BU1").Value = "sales_team.4.allocated_percentage [Float]"
$ws.Range("BV1").Value = "sales_team.4.incentives [Currency]"
$ws.Range("BW1").Value = "sales_team.5.sales_person [Link] [Sales Person]"
$ws.Range("BX1").Value = "sales_team.5.allocated_percentage [Float]"
$ws.Range("BY1").Value = "sales_team.5.incentives [Currency]"
$ws.Range("BZ1").Value = "default_sales_partner [Link] [Sales Partner]"
$ws.Range("CA1").Value = "default_commission_rate [Float]"
$ws.Range("CB1").Value = "so_required [Check]"
$ws.Range("CC1").Value = "dn_required [Check]"
$ws.Range("CD1").Value = "is_frozen [Check]"
$ws.Range("CE1").Value = "disabled [Check]"
$ws.Range("CF1").Value = "portal_users.1.user [Link] [User]"
$ws.Range("CG1").Value = "portal_users.2.user [Link] [User]"
$ws.Range("CH1").Value = "portal_users.3.user [Link] [User]"
$ws.Range("CI1").Value = "portal_users.4.user [Link] [User]"
$ws.Range("CJ1").Value = "portal_users.5.user [Link] [User]"
